# Commit: Update multi-assay metadata templates to use updated field descriptions
#
# This script:
#   1. Rewrites the 8 header-cell comments on the "SnareSeq2" sheet with revised
#      field descriptions (clearer wording + added examples).
#   2. Refreshes the "dataset_type" lookup sheet with the current HRAVS term list
#      (reordered + several additions/removals), growing it from 41 to 50 rows.
#   3. Re-points the dataset_type column validation at the larger lookup range.
#   4. Bumps the ".metadata" sheet's pav:createdOn timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the 8 header-cell comments on the "SnareSeq2" sheet.
# ---------------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("SnareSeq2")

$commentA1 = @"
(Required) The unique identifier from HuBMAP or SenNet for the sample (such as a
block, section, or suspension) used to perform the assay. For instance, in an
RNAseq assay, the parent sample would be the suspension, while in imaging
assays, it would be the tissue section. If the assay is derived from multiple
parent samples, this field should contain a comma-separated list of identifiers.
Example: HBM386.ZGKG.235, HBM672.MKPK.442
"@
[void]$wsMain.Range("A1").Comment.Text($commentA1)

$commentB1 = @"
A locally assigned identifier provided by the data provider for the dataset. It
is used to reference an external metadata record that may be maintained
independently, enabling traceability and supporting provenance tracking.
Example: Visium_9OLC_A4_S1
"@
[void]$wsMain.Range("B1").Comment.Text($commentB1)

$commentC1 = @"
(Required) The DOI for the protocols.io page that details the assay or the
procedures used for sample procurement and preparation. For example, in the case
of an imaging assay, the protocol may start with tissue section staining and end
with the generation of an OME-TIFF file. The documented protocol should also
include any image processing steps involved in producing the final OME-TIFF.
Example: https://dx.doi.org/10.17504/protocols.io.eq2lyno9qvx9/v1
"@
[void]$wsMain.Range("C1").Comment.Text($commentC1)

$commentD1 = @"
(Required) The specific type of dataset being produced. Example: RNAseq
"@
[void]$wsMain.Range("D1").Comment.Text($commentD1)

$commentE1 = @"
(Required) The name of the file containing the ORCID IDs for all contributors to
this dataset. Example: ./contributors.csv
"@
[void]$wsMain.Range("E1").Comment.Text($commentE1)

$commentF1 = @"
(Required) The top-level directory containing the raw and/or processed data. For
a single dataset upload, this might be represented as ".", whereas for a data
upload containing multiple datasets, this would be the directory name for the
respective dataset. For example, if the data is within a directory named
"TEST001-RK", use the syntax "./TEST001-RK" for this field. If there are
multiple directory levels, use the format "./TEST001-RK/Run1/Pass2", where
"Pass2" is the subdirectory where the single dataset's data is stored. This is
an internal metadata field used solely for data ingestion. Example: ./TEST001-RK
"@
[void]$wsMain.Range("F1").Comment.Text($commentF1)

$commentG1 = @"
The number of PCR cycles performed following the Chromium Controller step and
before the suspension is separated and library construction begins. Example: 7
"@
[void]$wsMain.Range("G1").Comment.Text($commentG1)

$commentH1 = @"
(Required) The unique string identifier for the metadata specification version,
which is easily interpretable by computers for purposes of data validation and
processing. Example: 22bc762a-5020-419d-b170-24253ed9e8d9
"@
[void]$wsMain.Range("H1").Comment.Text($commentH1)

# ---------------------------------------------------------------------------
# 2) Rebuild the "dataset_type" lookup sheet: refreshed list of (label, URL)
#    pairs, now 50 rows (was 41).
# ---------------------------------------------------------------------------
$wsTypes = $wb.Worksheets.Item("dataset_type")

$wsTypes.Range("A1").Value = "HiFi-Slide"
$wsTypes.Range("B1").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000195"
$wsTypes.Range("A2").Value = "SNARE-seq2"
$wsTypes.Range("B2").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000264"
$wsTypes.Range("A3").Value = "COMET"
$wsTypes.Range("B3").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000458"
$wsTypes.Range("A4").Value = "Visium (no probes)"
$wsTypes.Range("B4").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000302"
$wsTypes.Range("A5").Value = "DESI"
$wsTypes.Range("B5").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000204"
$wsTypes.Range("A6").Value = "Confocal"
$wsTypes.Range("B6").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000206"
$wsTypes.Range("A7").Value = "Stereo-seq"
$wsTypes.Range("B7").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000385"
$wsTypes.Range("A8").Value = "Visium (with probes)"
$wsTypes.Range("B8").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000303"
$wsTypes.Range("A9").Value = "Molecular Cartography"
$wsTypes.Range("B9").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000217"
$wsTypes.Range("A10").Value = "DBiT-seq"
$wsTypes.Range("B10").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000428"
$wsTypes.Range("A11").Value = "Seq-Scope"
$wsTypes.Range("B11").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000390"
$wsTypes.Range("A12").Value = "CosMx Transcriptomics"
$wsTypes.Range("B12").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000218"
$wsTypes.Range("A13").Value = "CyCIF"
$wsTypes.Range("B13").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000200"
$wsTypes.Range("A14").Value = "Light Sheet"
$wsTypes.Range("B14").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000168"
$wsTypes.Range("A15").Value = "seqFISH"
$wsTypes.Range("B15").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000397"
$wsTypes.Range("A16").Value = "ATACseq"
$wsTypes.Range("B16").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000309"
$wsTypes.Range("A17").Value = "CosMx Proteomics"
$wsTypes.Range("B17").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000435"
$wsTypes.Range("A18").Value = "Singular Genomics G4X"
$wsTypes.Range("B18").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000429"
$wsTypes.Range("A19").Value = "Visium HD"
$wsTypes.Range("B19").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000451"
$wsTypes.Range("A20").Value = "MERFISH"
$wsTypes.Range("B20").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000221"
$wsTypes.Range("A21").Value = "10X Multiome"
$wsTypes.Range("B21").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000215"
$wsTypes.Range("A22").Value = "4i"
$wsTypes.Range("B22").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000447"
$wsTypes.Range("A23").Value = "PhenoCycler"
$wsTypes.Range("B23").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000199"
$wsTypes.Range("A24").Value = "Second Harmonic Generation (SHG)"
$wsTypes.Range("B24").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000208"
$wsTypes.Range("A25").Value = "Thick section Multiphoton MxIF"
$wsTypes.Range("B25").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000207"
$wsTypes.Range("A26").Value = "CyTOF"
$wsTypes.Range("B26").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000407"
$wsTypes.Range("A27").Value = "Olink"
$wsTypes.Range("B27").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000441"
$wsTypes.Range("A28").Value = "MIBI"
$wsTypes.Range("B28").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000172"
$wsTypes.Range("A29").Value = "Auto-fluorescence"
$wsTypes.Range("B29").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000205"
$wsTypes.Range("A30").Value = "FACS"
$wsTypes.Range("B30").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000440"
$wsTypes.Range("A31").Value = "Xenium"
$wsTypes.Range("B31").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000219"
$wsTypes.Range("A32").Value = "SIMS"
$wsTypes.Range("B32").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000202"
$wsTypes.Range("A33").Value = "Cell DIVE"
$wsTypes.Range("B33").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000159"
$wsTypes.Range("A34").Value = "CODEX"
$wsTypes.Range("B34").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000160"
$wsTypes.Range("A35").Value = "GeoMx (NGS)"
$wsTypes.Range("B35").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000300"
$wsTypes.Range("A36").Value = "MUSIC"
$wsTypes.Range("B36").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000328"
$wsTypes.Range("A37").Value = "Pixel-seqV2"
$wsTypes.Range("B37").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000450"
$wsTypes.Range("A38").Value = "MALDI"
$wsTypes.Range("B38").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000201"
$wsTypes.Range("A39").Value = "2D Imaging Mass Cytometry"
$wsTypes.Range("B39").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000296"
$wsTypes.Range("A40").Value = "Histology"
$wsTypes.Range("B40").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000197"
$wsTypes.Range("A41").Value = "Enhanced Stimulated Raman Spectroscopy (SRS)"
$wsTypes.Range("B41").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000209"
$wsTypes.Range("A42").Value = "DART-FISH"
$wsTypes.Range("B42").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000396"
$wsTypes.Range("A43").Value = "Resolve"
$wsTypes.Range("B43").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000384"
$wsTypes.Range("A44").Value = "RNAseq"
$wsTypes.Range("B44").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000310"
$wsTypes.Range("A45").Value = "LC-MS"
$wsTypes.Range("B45").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000194"
$wsTypes.Range("A46").Value = "nanoSPLITS"
$wsTypes.Range("B46").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000312"
$wsTypes.Range("A47").Value = "GeoMx (nCounter)"
$wsTypes.Range("B47").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000301"
$wsTypes.Range("A48").Value = "RNAseq (with probes)"
$wsTypes.Range("B48").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000311"
$wsTypes.Range("A49").Value = "MS Lipidomics"
$wsTypes.Range("B49").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000405"
$wsTypes.Range("A50").Value = "MPLEx"
$wsTypes.Range("B50").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000448"

# ---------------------------------------------------------------------------
# 3) Point the D-column validation list at the now-larger dataset_type range.
# ---------------------------------------------------------------------------
$dValidation = $wsMain.Range("D2:D1001").Validation
$dValidation.Formula1 = "'dataset_type'!`$A`$1:`$A`$50"

# ---------------------------------------------------------------------------
# 4) Bump the ".metadata" sheet's pav:createdOn timestamp.
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item(".metadata")
$wsMeta.Range("C2").Value = "2025-10-21T13:32:46-07:00"

